# Update the R^2 / RMSE / U columns with the re-run notebook results
# and refresh the per-column background-gradient cell shading to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -5.583
$ws.Range("D2").Value = 0.9422
$ws.Range("E2").Value = 2.3053
$ws.Range("D2").Interior.Color = 0xF4FCF6
$ws.Range("E2").Interior.Color = 0xF3FBF5

$ws.Range("C3").Value = -2.0117
$ws.Range("D3").Value = 0.8015
$ws.Range("E3").Value = 1.697
$ws.Range("D3").Interior.Color = 0xAEE1B5
$ws.Range("E3").Interior.Color = 0x7BC77A

$ws.Range("C4").Value = -0.529
$ws.Range("D4").Value = 0.7043
$ws.Range("E4").Value = 1.4554
$ws.Range("D4").Interior.Color = 0x73C16E
$ws.Range("E4").Interior.Color = 0x549F36

$ws.Range("C5").Value = 0.5202
$ws.Range("D5").Value = 0.4748
$ws.Range("E5").Value = 1.0281
$ws.Range("D5").Interior.Color = 0x1B4400
$ws.Range("E5").Interior.Color = 0x1B4400

$ws.Range("C6").Value = 0.3059
$ws.Range("D6").Value = 0.632
$ws.Range("E6").Value = 1.4002
$ws.Range("D6").Interior.Color = 0x55A037
$ws.Range("E6").Interior.Color = 0x4C942C

$ws.Range("C7").Value = 0.0128
$ws.Range("D7").Value = 0.7428
$ws.Range("E7").Value = 1.8204
$ws.Range("D7").Interior.Color = 0x8AD08D
$ws.Range("E7").Interior.Color = 0x97D79C

$ws.Range("C8").Value = -0.2273
$ws.Range("D8").Value = 0.8346
$ws.Range("E8").Value = 2.0663
$ws.Range("D8").Interior.Color = 0xC3EACA
$ws.Range("E8").Interior.Color = 0xCEEED4

$ws.Range("C9").Value = -0.3744
$ws.Range("D9").Value = 0.8899
$ws.Range("E9").Value = 2.1959
$ws.Range("D9").Interior.Color = 0xE1F5E5
$ws.Range("E9").Interior.Color = 0xE5F7E9

$ws.Range("C10").Value = -0.5355
$ws.Range("D10").Value = 0.9475
$ws.Range("E10").Value = 2.3218
$ws.Range("D10").Interior.Color = 0xF5FCF7
$ws.Range("E10").Interior.Color = 0xF5FCF7
